$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 50: date must stay a literal text string (matching the inline
# strings used for every other "MM/DD/YYYY" date in column A), not get
# auto-converted into a date serial number by the COM Value setter.
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "08/15/2025"
$ws.Range("A50").ClearFormats()

$ws.Range("B50").Value = 550.1630000000005
$ws.Range("C50").Value = 0.09088215674263801
$ws.Range("D50").Value = 50
